$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cells flagged as changed in the diff (C2 and C4: 0 -> 1)
$ws.Range("C2").Value = 1
$ws.Range("C4").Value = 1

# Append the new rows (5-33)
$data = @(
    ,@("Вон-Пфеффер", "Иван", 0)
    ,@("Золотарев", "Ермолай", 1)
    ,@("Зорников", "Максим", 1)
    ,@("Казанков", "Захар", 1)
    ,@("Кубанова", "Валерия", 0)
    ,@("Ломей", "Егор", 1)
    ,@("Подыряка", "Евгений", 0)
    ,@("Чернобровкин", "Андрей", 1)
    ,@("Яницкая", "Ева", 0)
    ,@("Бондаренко", "Георгий", 1)
    ,@("Бурлаков", "Илья", 0)
    ,@("Вольная", "Полина", 1)
    ,@("Кучеров", "Фёдор", 0)
    ,@("Пясецкий", "Глеб", 1)
    ,@("Суханова", "Екатерина", 1)
    ,@("Чернов", "Александр", 1)
    ,@("Бродская", "Ольга", 1)
    ,@("Глазунов", "Владислав", 1)
    ,@("Зиновьев", "Владимир", 1)
    ,@("Измайлова", "Эмилия", 1)
    ,@("Конышев", "Мирон", 1)
    ,@("Лагунов", "Никита", 0)
    ,@("Пашнев", "Владимир", 1)
    ,@("Семенюк", "Артемий", 1)
    ,@("Соловьев", "Ярослав", 0)
    ,@("Та", "Хань Зуй", 0)
    ,@("Гладыш", "Андрей", 1)
    ,@("Пашнюк", "Марк", 0)
    ,@("Столярова", "Александра", 1)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 5 + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
